# Natmi following Dr Hou advice
# Rebuild the LR-pair table (Wnt9a -> Fzd4) to cover the full 3x3
# Sending-cluster x Target-cluster combination across ECs/FAPs/sCs,
# replacing the previous partial (2x2, missing ECs) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt9a"
$ws.Cells.Item(2, 3).Value = "Fzd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2454823333333333
$ws.Cells.Item(2, 8).Value = 0.736447
$ws.Cells.Item(2, 9).Value = 0.0314842275540206
$ws.Cells.Item(2, 10).Value = 0.0314842275540206
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 22.41709
$ws.Cells.Item(2, 14).Value = 67.25127000000001
$ws.Cells.Item(2, 15).Value = 0.3988455747018376
$ws.Cells.Item(2, 16).Value = 0.3988455747018376
$ws.Cells.Item(2, 17).Value = 5.502999559743333
$ws.Cells.Item(2, 18).Value = 49.52699603769
$ws.Cells.Item(2, 19).Value = 0.01255734483282678
$ws.Cells.Item(2, 20).Value = 0.01255734483282678

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt9a"
$ws.Cells.Item(3, 3).Value = "Fzd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2454823333333333
$ws.Cells.Item(3, 8).Value = 0.736447
$ws.Cells.Item(3, 9).Value = 0.0314842275540206
$ws.Cells.Item(3, 10).Value = 0.0314842275540206
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 16.78189033333333
$ws.Cells.Item(3, 14).Value = 50.345671
$ws.Cells.Item(3, 15).Value = 0.2985839238983091
$ws.Cells.Item(3, 16).Value = 0.2985839238983091
$ws.Cells.Item(3, 17).Value = 4.119657596770778
$ws.Cells.Item(3, 18).Value = 37.07691837093699
$ws.Cells.Item(3, 19).Value = 0.009400684203986733
$ws.Cells.Item(3, 20).Value = 0.009400684203986735

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt9a"
$ws.Cells.Item(4, 3).Value = "Fzd4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2454823333333333
$ws.Cells.Item(4, 8).Value = 0.736447
$ws.Cells.Item(4, 9).Value = 0.0314842275540206
$ws.Cells.Item(4, 10).Value = 0.0314842275540206
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 17.00595566666667
$ws.Cells.Item(4, 14).Value = 51.017867
$ws.Cells.Item(4, 15).Value = 0.3025705013998533
$ws.Cells.Item(4, 16).Value = 0.3025705013998533
$ws.Cells.Item(4, 17).Value = 4.174661677616555
$ws.Cells.Item(4, 18).Value = 37.571955098549
$ws.Cells.Item(4, 19).Value = 0.009526198517207089
$ws.Cells.Item(4, 20).Value = 0.009526198517207091

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Wnt9a"
$ws.Cells.Item(5, 3).Value = "Fzd4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 6.016075000000001
$ws.Cells.Item(5, 8).Value = 18.048225
$ws.Cells.Item(5, 9).Value = 0.7715890252063808
$ws.Cells.Item(5, 10).Value = 0.7715890252063808
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 22.41709
$ws.Cells.Item(5, 14).Value = 67.25127000000001
$ws.Cells.Item(5, 15).Value = 0.3988455747018376
$ws.Cells.Item(5, 16).Value = 0.3988455747018376
$ws.Cells.Item(5, 17).Value = 134.86289472175
$ws.Cells.Item(5, 18).Value = 1213.76605249575
$ws.Cells.Item(5, 19).Value = 0.3077448681920696
$ws.Cells.Item(5, 20).Value = 0.3077448681920696

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt9a"
$ws.Cells.Item(6, 3).Value = "Fzd4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6.016075000000001
$ws.Cells.Item(6, 8).Value = 18.048225
$ws.Cells.Item(6, 9).Value = 0.7715890252063808
$ws.Cells.Item(6, 10).Value = 0.7715890252063808
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 16.78189033333333
$ws.Cells.Item(6, 14).Value = 50.345671
$ws.Cells.Item(6, 15).Value = 0.2985839238983091
$ws.Cells.Item(6, 16).Value = 0.2985839238983091
$ws.Cells.Item(6, 17).Value = 100.9611108871083
$ws.Cells.Item(6, 18).Value = 908.649997983975
$ws.Cells.Item(6, 19).Value = 0.2303840787829925
$ws.Cells.Item(6, 20).Value = 0.2303840787829925

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt9a"
$ws.Cells.Item(7, 3).Value = "Fzd4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 6.016075000000001
$ws.Cells.Item(7, 8).Value = 18.048225
$ws.Cells.Item(7, 9).Value = 0.7715890252063808
$ws.Cells.Item(7, 10).Value = 0.7715890252063808
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 17.00595566666667
$ws.Cells.Item(7, 14).Value = 51.017867
$ws.Cells.Item(7, 15).Value = 0.3025705013998533
$ws.Cells.Item(7, 16).Value = 0.3025705013998533
$ws.Cells.Item(7, 17).Value = 102.3091047373417
$ws.Cells.Item(7, 18).Value = 920.781942636075
$ws.Cells.Item(7, 19).Value = 0.2334600782313187
$ws.Cells.Item(7, 20).Value = 0.2334600782313187

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Wnt9a"
$ws.Cells.Item(8, 3).Value = "Fzd4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.535436666666667
$ws.Cells.Item(8, 8).Value = 4.606310000000001
$ws.Cells.Item(8, 9).Value = 0.1969267472395986
$ws.Cells.Item(8, 10).Value = 0.1969267472395986
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 22.41709
$ws.Cells.Item(8, 14).Value = 67.25127000000001
$ws.Cells.Item(8, 15).Value = 0.3988455747018376
$ws.Cells.Item(8, 16).Value = 0.3988455747018376
$ws.Cells.Item(8, 17).Value = 34.42002194596667
$ws.Cells.Item(8, 18).Value = 309.7801975137
$ws.Cells.Item(8, 19).Value = 0.0785433616769412
$ws.Cells.Item(8, 20).Value = 0.07854336167694122

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Wnt9a"
$ws.Cells.Item(9, 3).Value = "Fzd4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.535436666666667
$ws.Cells.Item(9, 8).Value = 4.606310000000001
$ws.Cells.Item(9, 9).Value = 0.1969267472395986
$ws.Cells.Item(9, 10).Value = 0.1969267472395986
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 16.78189033333333
$ws.Cells.Item(9, 14).Value = 50.345671
$ws.Cells.Item(9, 15).Value = 0.2985839238983091
$ws.Cells.Item(9, 16).Value = 0.2985839238983091
$ws.Cells.Item(9, 17).Value = 25.76752975377889
$ws.Cells.Item(9, 18).Value = 231.90776778401
$ws.Cells.Item(9, 19).Value = 0.05879916091132985
$ws.Cells.Item(9, 20).Value = 0.05879916091132987

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Wnt9a"
$ws.Cells.Item(10, 3).Value = "Fzd4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.535436666666667
$ws.Cells.Item(10, 8).Value = 4.606310000000001
$ws.Cells.Item(10, 9).Value = 0.1969267472395986
$ws.Cells.Item(10, 10).Value = 0.1969267472395986
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 17.00595566666667
$ws.Cells.Item(10, 14).Value = 51.017867
$ws.Cells.Item(10, 15).Value = 0.3025705013998533
$ws.Cells.Item(10, 16).Value = 0.3025705013998533
$ws.Cells.Item(10, 17).Value = 26.11156788230778
$ws.Cells.Item(10, 18).Value = 235.00411094077
$ws.Cells.Item(10, 19).Value = 0.05958422465132751
$ws.Cells.Item(10, 20).Value = 0.05958422465132752

